$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.953.36'
$ws.Range("E2").Value = '  -0.53%  '
$ws.Range("D3").Value = '1.858.31'
$ws.Range("E3").Value = '  -0.95%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.98'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.48%  '
$ws.Range("E6").Value = '  -0.05%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5143'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.86%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3827'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.18%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08242'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.55%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.109'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.76%  '
$ws.Range("E11").Value = '  -0.23%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.174'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.59%  '
$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.49'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.85%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.861.82'
$ws.Range("E14").Value = '  -0.72%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.239'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.93%  '
$ws.Range("E16").Value = '  -0.09%  '
$ws.Range("E17").Value = '  -0.25%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '90.35'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.83%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06638'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.64'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.88%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.005'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.71%  '
$ws.Range("D23").Value = '27.989.77'
$ws.Range("E23").Value = '  -0.53%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.04'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.59%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.241'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.43%  '
$ws.Range("D26").Value = '2.073.54'
$ws.Range("E26").Value = '  -0.67%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.504'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.78%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '157.43'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.09%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '20.42'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.48%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.44'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.25%  '
$ws.Range("E31").Value = '  +1.27%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.026'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.09%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.875'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.593'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.49%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.409'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.14%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02411'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.92%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06500'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2179'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.25%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6537'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.53%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.192'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.19%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.969'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.51%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.207'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.14%  '
$ws.Range("E43").Value = '  -3.35%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6113'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.13%  '
$ws.Range("E45").Value = '  -2.10%  '
$ws.Range("E46").Value = '  -0.29%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.670'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.06%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.008'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.88%  '
$ws.Range("E49").Value = '  -1.74%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '120.68'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.56%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '78.19'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.70%  '
